$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, reusing the existing header style (bold,
# centered, bordered) by copying the H1 header cell's formatting, then
# overwriting the copied text with the new header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Fill columns I (I0) and J (IF) for each data row: I0 is always 1, IF
# mirrors the existing IP value (column H) for that row.
for ($r = 2; $r -le 27; $r++) {
    $ipVal = $ws.Cells.Item($r, 8).Value2  # column H = IP
    $ws.Cells.Item($r, 9).Value = 1          # column I = I0
    $ws.Cells.Item($r, 10).Value = $ipVal    # column J = IF (same as H)
}
